$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing column C header/values and add new column D (role/skill)
$ws.Range("C1").Value = "role"
$ws.Range("C2").Value = "qa"
$ws.Range("C3").Value = "dev"

$ws.Range("D1").Value = "skill"
$ws.Range("D2").Value = "uft"
$ws.Range("D3").Value = "selenium"

# Reflect the selection state captured in the saved file
$ws.Range("A4:XFD13").Select()
